$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 821; $r -le 1029; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "Congo, The Democratic Republic of the") {
        $cell.Value2 = "Democratic Republic of Congo"
    }
}
